$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Final data for rows 2-29: row index, A-column label, then B..J numeric values
$data = @(
    @(2, "CRP12", 1, 0, 2, 0, 3, 2, 2, 0, 1),
    @(3, "CRP15", 1, 0, 2, 0, 0, 6, 4, 0, 1),
    @(4, "CRP16", 0, 1, 2, 1, 1, 1, 2, 0, 1),
    @(5, "CRP17", 6, 16, 4, 6, 6, 13, 10, 6, 12),
    @(6, "CRP21", 3, 3, 4, 3, 3, 0, 7, 0, 0),
    @(7, "CRP22", 6, 2, 5, 0, 1, 4, 5, 6, 0),
    @(8, "CRP23", 8, 3, 4, 2, 0, 2, 8, 0, 0),
    @(9, "CRP24", 3, 4, 3, 9, 14, 15, 6, 12, 5),
    @(10, "CRP25", 4, 2, 9, 2, 5, 6, 7, 5, 8),
    @(11, "CRP26", 3, 3, 1, 0, 0, 10, 10, 1, 1),
    @(12, "CRP31", 8, 0, 10, 4, 2, 7, 10, 2, 12),
    @(13, "CRP32", 9, 9, 11, 9, 11, 17, 16, 10, 4),
    @(14, "CRP33", 3, 3, 7, 5, 4, 3, 9, 11, 4),
    @(15, "CRP34", 6, 6, 9, 7, 5, 6, 7, 1, 4),
    @(16, "CRP36", 12, 1, 8, 11, 9, 8, 9, 3, 3),
    @(17, "CRP37", 10, 4, 2, 1, 4, 10, 11, 5, 11),
    @(18, "CRP41", 3, 4, 2, 3, 0, 9, 7, 1, 4),
    @(19, "CRP42", 9, 10, 0, 17, 9, 18, 7, 0, 2),
    @(20, "CRP43", 11, 7, 5, 6, 5, 12, 17, 1, 0),
    @(21, "CRP44", 0, 0, 2, 0, 3, 1, 0, 0, 2),
    @(22, "CRP45", 5, 0, 8, 1, 1, 11, 9, 1, 1),
    @(23, "CRP46", 18, 4, 9, 11, 20, 16, 11, 5, 2),
    @(24, "CRP51", 10, 4, 5, 1, 1, 10, 11, 3, 5),
    @(25, "CRP52", 7, 0, 7, 5, 6, 6, 8, 1, 3),
    @(26, "CRP53", 11, 4, 8, 0, 6, 7, 8, 5, 9),
    @(27, "CRP54", 10, 2, 1, 12, 10, 10, 4, 3, 4),
    @(28, "CRP62", 2, 3, 2, 0, 0, 4, 2, 0, 1),
    @(29, "CRP63", 4, 4, 4, 5, 6, 4, 5, 6, 3)
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value2 = $row[1]
    for ($c = 2; $c -le 10; $c++) {
        $ws.Cells.Item($r, $c).Value2 = $row[$c]
    }
}

# Remove the now-obsolete trailing rows (previously CRP64 / CRP65)
$ws.Rows.Item(31).Delete()
$ws.Rows.Item(30).Delete()
